# TC21_Canine_Filter_Breed-Doberman.xlsx — "Fixed ICDC breed all testcases"
#
# The StatQuery column (C) on the "startup" sheet holds a long Cypher query
# string shared by cells C2:C4 (one shared string). Replace it with the new,
# simplified query, then update the row heights (which Excel auto-shrank
# because the new text is shorter) and restore the sheet's last-used
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Doberman Pinscher']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# StatQuery cells for CasesTab / SamplesTab / FilesTab rows.
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# The shorter query text re-wraps into fewer lines, so Excel's row autofit
# shrinks these rows down from the previous (capped) 409.6pt.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# Restore the sheet's saved selection/view state.
$ws.Range("B4:B5").Select()
